# Update shadow-model conversion values (D1/L1 = D2/L2 --> D2 = ...)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D4").Value = 8602.31615895128
$ws.Range("D5").Value = 8602.31615895128

$ws.Range("D9").Value = 11769.11175652173
$ws.Range("D10").Value = 11769.11175652173

$ws.Range("D14").Value = 11615.28384104874
$ws.Range("D15").Value = 11615.28384104874
